# "copy class from another done"
#
# Adds a new "DragNDropClass / Copy Class" feature row (row 8) to the
# "Praca s kliknutiami a klavesami" (Clickable) section of the sheet, with
# a new helper column F holding a footnote for that row, and moves the
# active selection to the new area (F13) the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 was a blank spacer row between the "Clickable" block (rows 1-7)
# and the "TableManager" block (row 9 onward); fill it in with the new
# feature instead of inserting a row, so everything below keeps its
# original row numbers.
$ws.Range("A8").Value2 = "DragNDropClass"
$ws.Range("B8").Value2 = "Copy Class"
$ws.Range("C8").Value2 = "skopiruje classu z jednej tabule do druhej"
$ws.Range("D8").Value2 = 'hold ,,F" + click na classu , potom hold ,,F" click na tabulu'
$ws.Range("E8").Value2 = "Clickable"
$ws.Range("F8").Value2 = "*neni su osetrene corner cases"

# New column F only carries data on row 8; give it the same boxed look as
# the rest of the table (thin right border) by cloning the format already
# used on the row instead of fabricating a fresh style entry.
$ws.Range("D8").Copy()
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("F8").Value2 = "*neni su osetrene corner cases"

# Give the new column a sensible custom width, matching the other
# description-style columns (C/D) in spirit.
$ws.Columns.Item(6).ColumnWidth = 30.8

# Author's final selection/cursor position ended up on the new column.
[void]$ws.Range("F13").Select()
